# Swap the data (all columns except A, which holds the running row id)
# between the given row pairs. This matches the diff where rows were
# re-sorted but the "id" column (A) stayed aligned to its row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $row1, $row2, $firstCol, $lastCol) {
    $r1 = $ws.Range("$firstCol$row1`:$lastCol$row1")
    $r2 = $ws.Range("$firstCol$row2`:$lastCol$row2")
    $tmp = $r1.Value2
    $r1.Value2 = $r2.Value2
    $r2.Value2 = $tmp
}

Swap-Rows $ws 105 106 "B" "AC"
Swap-Rows $ws 107 108 "B" "AC"
Swap-Rows $ws 132 133 "B" "AC"
